$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 138, pushing existing row 138 (blank separator) and the
# summary rows (sum [min] / sum [h] / sum [working weeks]) down by one.
$ws.Rows.Item(138).Insert()

# Update row 137: the "end" time (E137) changes from 0 to 00:15 (1.0416...E-2),
# which changes F137/G137 (shared formulas) automatically.
$ws.Range("E137").Value = 0.010416666666666666

# Fill in the new data row 138.
$ws.Range("A138").Value = 2014
$ws.Range("B138").Value = 7
$ws.Range("C138").Value = 13
$ws.Range("D138").Value = 0.33333333333333331
$ws.Range("E138").Value = 0.33333333333333331

# Extend the shared formulas in F and G down into the new row 138 (each row
# keeps the same relative formula the shared group already uses).
$ws.Range("F138").Formula = "=(E138-D138)*24*60"
$ws.Range("G138").Formula = "=F138/60"

# Apply the same number formats as the row above (time-spent columns) to the
# new data row, matching s="1" (hh:mm) for D/E, s="3" (0 decimals) for F and
# s="2" (2 decimals) for G.
$ws.Range("D138:E138").NumberFormat = $ws.Range("D137:E137").NumberFormat
$ws.Range("F138").NumberFormat = $ws.Range("F137").NumberFormat
$ws.Range("G138").NumberFormat = $ws.Range("G137").NumberFormat

# Update the SUM formula range (now row 138 is the last data row) and the
# dependent "sum [h]" / "sum [working weeks]" formulas, which got shifted to
# rows 140 and 141.
$ws.Range("F140").Formula = "=SUM(F2:F138)"
$ws.Range("F141").Formula = "=F140/60"
$ws.Range("F142").Formula = "=F141/38.5"
